# Update NATMI Fgf1-Fgfr4 LR-pair sheet with refreshed TPM-based statistics.
# Adds the "Inflammatory-Mac" target cluster (3 new rows) and refreshes
# all recalculated metric columns (G:T) for the existing sender/receiver pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.8775636666666666
$ws.Range("H2").Value = 2.632691
$ws.Range("I2").Value = 0.1887436506618166
$ws.Range("J2").Value = 0.2083714858314108
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.363908
$ws.Range("N2").Value = 1.091724
$ws.Range("O2").Value = 0.01118972054940699
$ws.Range("P2").Value = 0.01663265480083782
$ws.Range("Q2").Value = 0.3193524388093333
$ws.Range("R2").Value = 2.874171949284
$ws.Range("S2").Value = 0.002111988706380624
$ws.Range("T2").Value = 0.003465770994171524

# Row 3
$ws.Range("G3").Value = 0.8775636666666666
$ws.Range("H3").Value = 2.632691
$ws.Range("I3").Value = 0.1887436506618166
$ws.Range("J3").Value = 0.2083714858314108
$ws.Range("N3").Value = 0.460698
$ws.Range("O3").Value = 0.004721964413781051
$ws.Range("P3").Value = 0.007018835164781924
$ws.Range("Q3").Value = 0.1347639420353333
$ws.Range("R3").Value = 1.212875478318
$ws.Range("S3").Value = 0.0008912408017522201
$ws.Range("T3").Value = 0.001462525112091365

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.8775636666666666
$ws.Range("H4").Value = 2.632691
$ws.Range("I4").Value = 0.1887436506618166
$ws.Range("J4").Value = 0.2083714858314108
$ws.Range("K4").Value = 1.0
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06252866666666666
$ws.Range("N4").Value = 0.187586
$ws.Range("O4").Value = 0.001922679101110775
$ws.Range("P4").Value = 0.002857913889838424
$ws.Range("Q4").Value = 0.05487288599177777
$ws.Range("R4").Value = 0.493855973926
$ws.Range("S4").Value = 0.0003628934725948277
$ws.Range("T4").Value = 0.0005955077636038592

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.8775636666666666
$ws.Range("H5").Value = 2.632691
$ws.Range("I5").Value = 0.1887436506618166
$ws.Range("J5").Value = 0.2083714858314108
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 31.927516
$ws.Range("N5").Value = 63.85503199999999
$ws.Range("O5").Value = 0.9817315966582778
$ws.Range("P5").Value = 0.9728454303033116
$ws.Range("Q5").Value = 28.01842800851866
$ws.Range("R5").Value = 168.110568051112
$ws.Range("S5").Value = 0.1852956055233374
$ws.Range("T5").Value = 0.2027132477965992

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 0.8775636666666666
$ws.Range("H6").Value = 2.632691
$ws.Range("I6").Value = 0.1887436506618166
$ws.Range("J6").Value = 0.2083714858314108
$ws.Range("K6").Value = 1.0
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01411566666666667
$ws.Range("N6").Value = 0.042347
$ws.Range("O6").Value = 0.0004340392774233579
$ws.Range("P6").Value = 0.0006451658412300904
$ws.Range("Q6").Value = 0.01238739619744444
$ws.Range("R6").Value = 0.111486565777
$ws.Range("S6").Value = 0.00008192215775150155
$ws.Range("T6").Value = 0.000134434164944786

# Row 7
$ws.Range("D7").Value = "ECs"
$ws.Range("G7").Value = 2.458038666666667
$ws.Range("H7").Value = 7.374116000000001
$ws.Range("I7").Value = 0.5286672739959656
$ws.Range("J7").Value = 0.5836444564186148
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 0.363908
$ws.Range("N7").Value = 1.091724
$ws.Range("O7").Value = 0.01118972054940699
$ws.Range("P7").Value = 0.01663265480083782
$ws.Range("Q7").Value = 0.8944999351093333
$ws.Range("R7").Value = 8.050499415984
$ws.Range("S7").Value = 0.005915639059631633
$ws.Range("T7").Value = 0.009707556770033455

# Row 8
$ws.Range("D8").Value = "FAPs"
$ws.Range("G8").Value = 2.458038666666667
$ws.Range("H8").Value = 7.374116000000001
$ws.Range("I8").Value = 0.5286672739959656
$ws.Range("J8").Value = 0.5836444564186148
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.153566
$ws.Range("N8").Value = 0.460698
$ws.Range("O8").Value = 0.004721964413781051
$ws.Range("P8").Value = 0.007018835164781924
$ws.Range("Q8").Value = 0.3774711658853334
$ws.Range("R8").Value = 3.397240492968
$ws.Range("S8").Value = 0.002496348054539586
$ws.Range("T8").Value = 0.004096504234441006

# Row 9
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("G9").Value = 2.458038666666667
$ws.Range("H9").Value = 7.374116000000001
$ws.Range("I9").Value = 0.5286672739959656
$ws.Range("J9").Value = 0.5836444564186148
$ws.Range("M9").Value = 0.06252866666666666
$ws.Range("N9").Value = 0.187586
$ws.Range("O9").Value = 0.001922679101110775
$ws.Range("P9").Value = 0.002857913889838424
$ws.Range("Q9").Value = 0.1536978804417778
$ws.Range("R9").Value = 1.383280923976
$ws.Range("S9").Value = 0.001016457519153247
$ws.Range("T9").Value = 0.001668005598725956

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3.0
$ws.Range("G10").Value = 2.458038666666667
$ws.Range("H10").Value = 7.374116000000001
$ws.Range("I10").Value = 0.5286672739959656
$ws.Range("J10").Value = 0.5836444564186148
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 31.927516
$ws.Range("N10").Value = 63.85503199999999
$ws.Range("O10").Value = 0.9817315966582778
$ws.Range("P10").Value = 0.9728454303033116
$ws.Range("Q10").Value = 78.47906885861867
$ws.Range("R10").Value = 470.874413151712
$ws.Range("S10").Value = 0.5190093670010385
$ws.Range("T10").Value = 0.5677958423487097

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3.0
$ws.Range("G11").Value = 2.458038666666667
$ws.Range("H11").Value = 7.374116000000001
$ws.Range("I11").Value = 0.5286672739959656
$ws.Range("J11").Value = 0.5836444564186148
$ws.Range("M11").Value = 0.01411566666666667
$ws.Range("N11").Value = 0.042347
$ws.Range("O11").Value = 0.0004340392774233579
$ws.Range("P11").Value = 0.0006451658412300904
$ws.Range("Q11").Value = 0.03469685447244445
$ws.Range("R11").Value = 0.3122716902520001
$ws.Range("S11").Value = 0.0002294623616025853
$ws.Range("T11").Value = 0.0003765474667045945

# Row 12
$ws.Range("D12").Value = "ECs"
$ws.Range("G12").Value = 1.313898
$ws.Range("H12").Value = 2.627796
$ws.Range("I12").Value = 0.2825890753422177
$ws.Range("J12").Value = 0.2079840577499744
$ws.Range("K12").Value = 3.0
$ws.Range("M12").Value = 0.363908
$ws.Range("N12").Value = 1.091724
$ws.Range("O12").Value = 0.01118972054940699
$ws.Range("P12").Value = 0.01663265480083782
$ws.Range("Q12").Value = 0.4781379933839999
$ws.Range("R12").Value = 2.868827960304
$ws.Range("S12").Value = 0.003162092783394735
$ws.Range("T12").Value = 0.003459327036632842

# Row 13
$ws.Range("D13").Value = "FAPs"
$ws.Range("G13").Value = 1.313898
$ws.Range("H13").Value = 2.627796
$ws.Range("I13").Value = 0.2825890753422177
$ws.Range("J13").Value = 0.2079840577499744
$ws.Range("M13").Value = 0.153566
$ws.Range("N13").Value = 0.460698
$ws.Range("O13").Value = 0.004721964413781051
$ws.Range("P13").Value = 0.007018835164781924
$ws.Range("Q13").Value = 0.201770060268
$ws.Range("R13").Value = 1.210620361608
$ws.Range("S13").Value = 0.001334375557489244
$ws.Range("T13").Value = 0.001459805818249555

# Row 14  (new row)
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Fgf1"
$ws.Range("C14").Value = "Fgfr4"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 1.313898
$ws.Range("H14").Value = 2.627796
$ws.Range("I14").Value = 0.2825890753422177
$ws.Range("J14").Value = 0.2079840577499744
$ws.Range("K14").Value = 1.0
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.06252866666666666
$ws.Range("N14").Value = 0.187586
$ws.Range("O14").Value = 0.001922679101110775
$ws.Range("P14").Value = 0.002857913889838424
$ws.Range("Q14").Value = 0.08215629007599999
$ws.Range("R14").Value = 0.492937740456
$ws.Range("S14").Value = 0.0005433281093627004
$ws.Range("T14").Value = 0.0005944005275086088

# Row 15  (new row)
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Fgf1"
$ws.Range("C15").Value = "Fgfr4"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 2.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 1.313898
$ws.Range("H15").Value = 2.627796
$ws.Range("I15").Value = 0.2825890753422177
$ws.Range("J15").Value = 0.2079840577499744
$ws.Range("K15").Value = 2.0
$ws.Range("L15").Value = 1.0
$ws.Range("M15").Value = 31.927516
$ws.Range("N15").Value = 63.85503199999999
$ws.Range("O15").Value = 0.9817315966582778
$ws.Range("P15").Value = 0.9728454303033116
$ws.Range("Q15").Value = 41.949499417368
$ws.Range("R15").Value = 167.797997669472
$ws.Range("S15").Value = 0.2774266241339018
$ws.Range("T15").Value = 0.2023363401580027

# Row 16  (new row)
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Fgf1"
$ws.Range("C16").Value = "Fgfr4"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 1.313898
$ws.Range("H16").Value = 2.627796
$ws.Range("I16").Value = 0.2825890753422177
$ws.Range("J16").Value = 0.2079840577499744
$ws.Range("K16").Value = 1.0
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.01411566666666667
$ws.Range("N16").Value = 0.042347
$ws.Range("O16").Value = 0.0004340392774233579
$ws.Range("P16").Value = 0.0006451658412300904
$ws.Range("Q16").Value = 0.018546546202
$ws.Range("R16").Value = 0.111279277212
$ws.Range("S16").Value = 0.000122654758069271
$ws.Range("T16").Value = 0.0001341842095807099

